$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Data")

# New flight arrival entry to append after the last existing row (row 28 -> row 29)
$row = 29

$ws.Cells.Item($row, 1).Value = 28
$ws.Cells.Item($row, 2).Value = "Monday, Jan 16"
$ws.Cells.Item($row, 3).Value = "8:00 AM"
$ws.Cells.Item($row, 4).Value = "FR6111"
$ws.Cells.Item($row, 5).Value = "Gdansk"
$ws.Cells.Item($row, 6).Value = "(GDN)"
$ws.Cells.Item($row, 7).Value = "Ryanair "
$ws.Cells.Item($row, 8).Value = "B738"
$ws.Cells.Item($row, 9).Value = "(SP-RSW)"
$ws.Cells.Item($row, 10).Value = "7:46 AM"
$ws.Cells.Item($row, 12).Value = "0 hours, -14 minutes"
